$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04339299999999999
$ws.Range("H2").Value = 0.130179
$ws.Range("I2").Value = 0.0698021577815419
$ws.Range("J2").Value = 0.0698021577815419
$ws.Range("M2").Value = 0.3883076666666667
$ws.Range("N2").Value = 1.164923
$ws.Range("O2").Value = 0.1188638477168776
$ws.Range("P2").Value = 0.1188638477168776
$ws.Range("Q2").Value = 0.01684983457966666
$ws.Range("R2").Value = 0.151648511217
$ws.Range("S2").Value = 0.008296953052854661
$ws.Range("T2").Value = 0.008296953052854661
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04339299999999999
$ws.Range("H3").Value = 0.130179
$ws.Range("I3").Value = 0.0698021577815419
$ws.Range("J3").Value = 0.0698021577815419
$ws.Range("O3").Value = 0.6829215134520935
$ws.Range("P3").Value = 0.6829215134520935
$ws.Range("Q3").Value = 0.09680920442666664
$ws.Range("R3").Value = 0.8712828398399999
$ws.Range("S3").Value = 0.04766939523439243
$ws.Range("T3").Value = 0.04766939523439243
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04339299999999999
$ws.Range("H4").Value = 0.130179
$ws.Range("I4").Value = 0.0698021577815419
$ws.Range("J4").Value = 0.0698021577815419
$ws.Range("O4").Value = 0.1982146388310289
$ws.Range("P4").Value = 0.1982146388310289
$ws.Range("Q4").Value = 0.028098399469
$ws.Range("R4").Value = 0.252885595221
$ws.Range("S4").Value = 0.01383580949429482
$ws.Range("T4").Value = 0.01383580949429482
$ws.Range("I5").Value = 0.6764796878879081
$ws.Range("J5").Value = 0.6764796878879081
$ws.Range("M5").Value = 0.3883076666666667
$ws.Range("N5").Value = 1.164923
$ws.Range("O5").Value = 0.1188638477168776
$ws.Range("P5").Value = 0.1188638477168776
$ws.Range("Q5").Value = 0.1632982589605555
$ws.Range("R5").Value = 1.469684330645
$ws.Range("S5").Value = 0.08040897860466921
$ws.Range("T5").Value = 0.08040897860466921
$ws.Range("I6").Value = 0.6764796878879081
$ws.Range("J6").Value = 0.6764796878879081
$ws.Range("O6").Value = 0.6829215134520935
$ws.Range("P6").Value = 0.6829215134520935
$ws.Range("S6").Value = 0.4619825322720101
$ws.Range("T6").Value = 0.4619825322720101
$ws.Range("I7").Value = 0.6764796878879081
$ws.Range("J7").Value = 0.6764796878879081
$ws.Range("O7").Value = 0.1982146388310289
$ws.Range("P7").Value = 0.1982146388310289
$ws.Range("S7").Value = 0.1340881770112289
$ws.Range("T7").Value = 0.1340881770112289
$ws.Range("I8").Value = 0.2537181543305499
$ws.Range("J8").Value = 0.2537181543305499
$ws.Range("M8").Value = 0.3883076666666667
$ws.Range("N8").Value = 1.164923
$ws.Range("O8").Value = 0.1188638477168776
$ws.Range("P8").Value = 0.1188638477168776
$ws.Range("Q8").Value = 0.06124608559677779
$ws.Range("R8").Value = 0.551214770371
$ws.Range("S8").Value = 0.03015791605935374
$ws.Range("T8").Value = 0.03015791605935374
$ws.Range("I9").Value = 0.2537181543305499
$ws.Range("J9").Value = 0.2537181543305499
$ws.Range("O9").Value = 0.6829215134520935
$ws.Range("P9").Value = 0.6829215134520935
$ws.Range("S9").Value = 0.173269585945691
$ws.Range("T9").Value = 0.173269585945691
$ws.Range("I10").Value = 0.2537181543305499
$ws.Range("J10").Value = 0.2537181543305499
$ws.Range("O10").Value = 0.1982146388310289
$ws.Range("P10").Value = 0.1982146388310289
$ws.Range("S10").Value = 0.05029065232550522
$ws.Range("T10").Value = 0.05029065232550522
